$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.580.35"
$ws.Range("E2").Value = "  -2.93%  "
$ws.Range("D3").Value = "3.322.72"
$ws.Range("E3").Value = "  -4.80%  "
$r = $ws.Range("D4")
$r.Value = "'0.999"
$r.Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "
$r = $ws.Range("D5")
$r.Value = "'550.68"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "
$r = $ws.Range("D6")
$r.Value = "'173.55"
$r.Style = "Normal"
$ws.Range("E6").Value = "  -3.24%  "
$r = $ws.Range("D7")
$r.Value = "'0.611"
$r.Style = "Normal"
$ws.Range("E7").Value = "  -4.71%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "3.315.17"
$ws.Range("E9").Value = "  -4.88%  "
$r = $ws.Range("D10")
$r.Value = "'0.613"
$r.Style = "Normal"
$ws.Range("E10").Value = "  -3.04%  "
$ws.Range("E11").Value = "  -1.78%  "
$r = $ws.Range("D12")
$r.Value = "'53.20"
$r.Style = "Normal"
$ws.Range("E12").Value = "  -1.33%  "
$r = $ws.Range("D13")
$r.Value = "'0.0000265"
$r.Style = "Normal"
$ws.Range("E13").Value = "  -2.71%  "
$r = $ws.Range("D14")
$r.Value = "'8.91"
$r.Style = "Normal"
$ws.Range("E14").Value = "  -2.88%  "
$ws.Range("D15").Value = "3.848.69"
$ws.Range("E15").Value = "  -4.95%  "
$r = $ws.Range("D16")
$r.Value = "'18.24"
$r.Style = "Normal"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("E17").Value = "  -3.26%  "
$ws.Range("D18").Value = "3.316.03"
$ws.Range("E18").Value = "  -5.01%  "
$r = $ws.Range("D19")
$r.Value = "'11.72"
$r.Style = "Normal"
$ws.Range("E19").Value = "  -3.98%  "
$ws.Range("D20").Value = "63.433.85"
$ws.Range("E20").Value = "  -3.20%  "
$r = $ws.Range("D21")
$r.Value = "'0.969"
$r.Style = "Normal"
$ws.Range("E21").Value = "  -2.59%  "
$r = $ws.Range("D22")
$r.Value = "'430.09"
$r.Style = "Normal"
$ws.Range("E22").Value = "  +3.81%  "
$r = $ws.Range("D23")
$r.Value = "'4.44"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +8.15%  "
$r = $ws.Range("D24")
$r.Value = "'4.06"
$r.Style = "Normal"
$ws.Range("E24").Value = "  +0.11%  "
$r = $ws.Range("D25")
$r.Value = "'13.31"
$r.Style = "Normal"
$ws.Range("E25").Value = "  +4.29%  "
$r = $ws.Range("D26")
$r.Value = "'83.26"
$r.Style = "Normal"
$ws.Range("E26").Value = "  -2.86%  "
$r = $ws.Range("D27")
$r.Value = "'10.64"
$r.Style = "Normal"
$ws.Range("E27").Value = "  -1.19%  "
$r = $ws.Range("D28")
$r.Value = "'2.73"
$r.Style = "Normal"
$ws.Range("E28").Value = "  -4.29%  "
$r = $ws.Range("D29")
$r.Value = "'8.69"
$r.Style = "Normal"
$ws.Range("E29").Value = "  -3.81%  "
$r = $ws.Range("D30")
$r.Value = "'29.26"
$r.Style = "Normal"
$ws.Range("E30").Value = "  -3.58%  "
$r = $ws.Range("D31")
$r.Value = "'6.46"
$r.Style = "Normal"
$ws.Range("E31").Value = "  +0.17%  "
$r = $ws.Range("D32")
$r.Value = "'11.38"
$r.Style = "Normal"
$ws.Range("E32").Value = "  -2.16%  "
$r = $ws.Range("D33")
$r.Value = "'578.57"
$r.Style = "Normal"
$ws.Range("E33").Value = "  -6.00%  "
$ws.Range("E34").Value = "  -3.45%  "
$r = $ws.Range("D35")
$r.Value = "'58.20"
$r.Style = "Normal"
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("E37").Value = "  -1.08%  "
$r = $ws.Range("D38")
$r.Value = "'3.46"
$r.Style = "Normal"
$ws.Range("E38").Value = "  +5.56%  "
$r = $ws.Range("D39")
$r.Value = "'35.15"
$r.Style = "Normal"
$ws.Range("E39").Value = "  -5.28%  "
$ws.Range("D40").Value = "0.0₃0740"
$ws.Range("E40").Value = "  -6.41%  "
$r = $ws.Range("D41")
$r.Value = "'0.365"
$r.Style = "Normal"
$ws.Range("E41").Value = "  -4.02%  "
$ws.Range("D42").Value = "3.121.92"
$ws.Range("E42").Value = "  -7.50%  "
$r = $ws.Range("D43")
$r.Value = "'0.997"
$r.Style = "Normal"
$ws.Range("E43").Value = "  -0.25%  "
$r = $ws.Range("D44")
$r.Value = "'2.80"
$r.Style = "Normal"
$ws.Range("E44").Value = "  -1.77%  "
$r = $ws.Range("D45")
$r.Value = "'3.17"
$r.Style = "Normal"
$ws.Range("E45").Value = "  -3.39%  "
$r = $ws.Range("D46")
$r.Value = "'0.0403"
$r.Style = "Normal"
$ws.Range("E46").Value = "  -3.17%  "
$ws.Range("E47").Value = "  -3.85%  "
$r = $ws.Range("D48")
$r.Value = "'2.60"
$r.Style = "Normal"
$ws.Range("E48").Value = "  -6.23%  "
$ws.Range("E49").Value = "  -4.04%  "
$r = $ws.Range("D50")
$r.Value = "'132.94"
$r.Style = "Normal"
$ws.Range("E50").Value = "  -3.21%  "
$r = $ws.Range("D51")
$r.Value = "'8.06"
$r.Style = "Normal"
$ws.Range("E51").Value = "  -4.28%  "
